# Update build/version timestamps for release "mines - January 30"
# Old build timestamp: January 30 2026 16.19.47 EST
# New build timestamp: February 02 2026 12.49.33 EST

$wb = $excel.ActiveWorkbook

$oldStamp = "January 30 2026 16.19.47 EST"
$newStamp = "February 02 2026 12.49.33 EST"

$aboutWs = $wb.Worksheets.Item("About")
$dataWs = $wb.Worksheets.Item("Boundaries and methane sources")

# A2: "Version: mines - January 30 (built on January 30 2026 16.19.47 EST)"
$a2 = $aboutWs.Range("A2").Value()
$aboutWs.Range("A2").Value = $a2.Replace($oldStamp, $newStamp)

# A6: Recommended Citation text with the same build stamp embedded
$a6 = $aboutWs.Range("A6").Value()
$aboutWs.Range("A6").Value = $a6.Replace($oldStamp, $newStamp)

# S2:S11 on the data sheet: "mines - January 30 (built on January 30 2026 16.19.47 EST)"
for ($row = 2; $row -le 11; $row++) {
    $cell = $dataWs.Cells.Item($row, 19)
    $val = $cell.Value()
    $cell.Value = $val.Replace($oldStamp, $newStamp)
}
